# Generate Report for Handback
#
# The localization-status report tracks two source files
# (37ea1dd4-... and c003bd65-...). The c003bd65 file has just been
# handed back (it is now "in sync with en-US"), so:
#   - it moves to the top of the Overview summary sheet
#   - its row in each locale sheet (zh-cn, de-de) gets a Status of
#     "Handed back: in sync with en-US" plus newly-populated
#     "Latest Target File" / "Latest Handback File" / "Latest Handback
#     DateTime" columns.
# The other file (37ea1dd4-...) is untouched content-wise but drops to
# the second row.

$wb = $excel.ActiveWorkbook

$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/{0}/e2e/{1}"
$xlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/{0}/ol-handoff/OpenLocalizationTestOrg/oltest.{1}/ci/ht/{2}"

$c003Commit = "43f7a690a20b905f31a8c3fdf488167a3321d2e8"
$c37eaCommit = "a769f066dcdfc66e2f1210d9ce9ee413c8966878"

$c003Md  = "c003bd65-8677-4b9d-aad3-abac071d090b.md"
$c37eaMd = "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.md"

$c003ZhXlf  = "c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.zh-cn.xlf"
$c003DeXlf  = "c003bd65-8677-4b9d-aad3-abac071d090b.d9051e598847c7ea9d5cb7a0011e8a68085e1619.de-de.xlf"
$c37eaZhXlf = "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.5e7a37b2da99f014721c2b7a355e1b9a70c6751c.zh-cn.xlf"
$c37eaDeXlf = "37ea1dd4-59ab-4d3f-9c37-96cb15cebf1c.5e7a37b2da99f014721c2b7a355e1b9a70c6751c.de-de.xlf"

$c003ZhCommit  = "cd824585b09f8bf94fa4886d3450c9c9e3636bd8"
$c003DeCommit  = "b1cb137ed09fbabdeb80df581d021d696143428d"
$c37eaZhCommit = "406abe3598a4c180756b3ecd0af98c86adb5b31d"
$c37eaDeCommit = "ff3a514f657f04f247309fad5bfb5fa8e767cbdd"

$statusHandedBack = "Handed back: in sync with en-US"
$statusReady       = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet: c003bd65 row moves to row 2, 37ea1dd4 drops to row 3.
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")

$ovw.Hyperlinks.Delete()

$ovw.Range("A2").Value = $c003Md
$ovw.Range("B2").Value = $statusHandedBack
$ovw.Range("C2").Value = $statusHandedBack
$ovw.Range("D2").Value = "2016-25-11 18:25:31"

$ovw.Range("A3").Value = $c37eaMd
$ovw.Range("B3").Value = $statusReady
$ovw.Range("C3").Value = $statusReady
$ovw.Range("D3").Value = "2016-25-11 18:25:12"

$ovw.Hyperlinks.Add($ovw.Range("A2"), ($mdUrl -f $c003Commit, $c003Md), "", "", $c003Md)
$ovw.Hyperlinks.Add($ovw.Range("A3"), ($mdUrl -f $c37eaCommit, $c37eaMd), "", "", $c37eaMd)

# ---------------------------------------------------------------------
# zh-cn sheet: c003bd65 row moves to row 2 (now with handback info),
# 37ea1dd4 drops to row 3.
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Hyperlinks.Delete()

$zh.Range("A2").Value = $c003Md
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = $statusHandedBack
$zh.Range("D2").Value = $c003ZhXlf
$zh.Range("E2").Value = "2016-03-11 18:25:28"
$zh.Range("F2").Value = $c003Md
$zh.Range("G2").Value = $c003ZhXlf
$zh.Range("H2").Value = "2016-03-11 18:25:46"
$zh.Range("I2").Value = "Include"

$zh.Range("A3").Value = $c37eaMd
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = $statusReady
$zh.Range("D3").Value = $c37eaZhXlf
$zh.Range("E3").Value = "2016-03-11 18:25:09"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("I3").Value = "Include"

$zh.Hyperlinks.Add($zh.Range("A2"), ($mdUrl -f $c003Commit, $c003Md), "", "", $c003Md)
$zh.Hyperlinks.Add($zh.Range("B2"), ($mdUrl -f $c003Commit, $c003Md), "", "", ".md")
$zh.Hyperlinks.Add($zh.Range("D2"), ($xlfUrl -f $c003ZhCommit, "zh-cn", $c003ZhXlf), "", "", $c003ZhXlf)
$zh.Hyperlinks.Add($zh.Range("F2"), ($mdUrl -f $c003Commit, $c003Md), "", "", $c003Md)
$zh.Hyperlinks.Add($zh.Range("G2"), ($xlfUrl -f $c003ZhCommit, "zh-cn", $c003ZhXlf), "", "", $c003ZhXlf)

$zh.Hyperlinks.Add($zh.Range("A3"), ($mdUrl -f $c37eaCommit, $c37eaMd), "", "", $c37eaMd)
$zh.Hyperlinks.Add($zh.Range("B3"), ($mdUrl -f $c37eaCommit, $c37eaMd), "", "", ".md")
$zh.Hyperlinks.Add($zh.Range("D3"), ($xlfUrl -f $c37eaZhCommit, "zh-cn", $c37eaZhXlf), "", "", $c37eaZhXlf)

# ---------------------------------------------------------------------
# de-de sheet: same shape as zh-cn.
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Hyperlinks.Delete()

$de.Range("A2").Value = $c003Md
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = $statusHandedBack
$de.Range("D2").Value = $c003DeXlf
$de.Range("E2").Value = "2016-03-11 18:25:31"
$de.Range("F2").Value = $c003Md
$de.Range("G2").Value = $c003DeXlf
$de.Range("H2").Value = "2016-03-11 18:25:54"
$de.Range("I2").Value = "Include"

$de.Range("A3").Value = $c37eaMd
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = $statusReady
$de.Range("D3").Value = $c37eaDeXlf
$de.Range("E3").Value = "2016-03-11 18:25:12"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("I3").Value = "Include"

$de.Hyperlinks.Add($de.Range("A2"), ($mdUrl -f $c003Commit, $c003Md), "", "", $c003Md)
$de.Hyperlinks.Add($de.Range("B2"), ($mdUrl -f $c003Commit, $c003Md), "", "", ".md")
$de.Hyperlinks.Add($de.Range("D2"), ($xlfUrl -f $c003DeCommit, "de-de", $c003DeXlf), "", "", $c003DeXlf)
$de.Hyperlinks.Add($de.Range("F2"), ($mdUrl -f $c003Commit, $c003Md), "", "", $c003Md)
$de.Hyperlinks.Add($de.Range("G2"), ($xlfUrl -f $c003DeCommit, "de-de", $c003DeXlf), "", "", $c003DeXlf)

$de.Hyperlinks.Add($de.Range("A3"), ($mdUrl -f $c37eaCommit, $c37eaMd), "", "", $c37eaMd)
$de.Hyperlinks.Add($de.Range("B3"), ($mdUrl -f $c37eaCommit, $c37eaMd), "", "", ".md")
$de.Hyperlinks.Add($de.Range("D3"), ($xlfUrl -f $c37eaDeCommit, "de-de", $c37eaDeXlf), "", "", $c37eaDeXlf)

Write-Host "Handback report regenerated."
